# "Add minor slide tweaks"
#
# The author edited the title slide (slide 1): the subtitle line read
# "[10] Cyber Law" and was bumped to "[11] Cyber Law" (the lecture's
# week/number tag was incremented, the rest of the line left alone).
#
# We locate the title placeholder shape on slide 1 and replace just the
# "[10] " prefix with "[11] " in place, leaving the surrounding run
# ("CITS1003 Introduction to Cybersecurity", the line break, and the
# trailing "Cyber Law" text) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleShape = $s.Shapes.Item("Title 1")
$tr = $titleShape.TextFrame.TextRange

$fullText = $tr.Text
$oldTag = "[10] "
$newTag = "[11] "

$tagStart = $fullText.IndexOf($oldTag)
if ($tagStart -ge 0) {
    # TextRange.Characters uses 1-based character positions.
    $tagRange = $tr.Characters($tagStart + 1, $oldTag.Length)
    $tagRange.Text = $newTag
}
